$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$fmt = $ws2.Range("AE2").NumberFormat
Write-Output "Format: $fmt"
$ws2.Range("AF2").NumberFormat = $fmt
$ws2.Range("AF2").Value = 43951
Write-Output "Done"
